$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3: Base Load / Market share_class_min -- drop the 2000 (M) value and
# raise every ramp year from 0.999 to 1.
# ---------------------------------------------------------------------------
$ws.Range("M3").ClearContents()
$ws.Range("N3:W3").Value = 1

# ---------------------------------------------------------------------------
# Row 17: Shoulder Load / Market share_class_min -- drop 2000 (M), new ramp.
# ---------------------------------------------------------------------------
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = 0.8
$ws.Range("O17").Value = 0.82499999999999996
$ws.Range("P17").Value = 0.85
$ws.Range("Q17").Value = 0.875
$ws.Range("R17").Value = 0.9
$ws.Range("S17:W17").Value = 1

# ---------------------------------------------------------------------------
# Row 24: Peak Load / Market share_class_min -- drop 2000 (M), new ramp.
# ---------------------------------------------------------------------------
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = 0.6
$ws.Range("O24").Value = 0.65
$ws.Range("P24").Value = 0.7
$ws.Range("Q24").Value = 0.75
$ws.Range("R24").Value = 0.8
$ws.Range("S24:W24").Value = 1

# ---------------------------------------------------------------------------
# New rows 27-33: Base Load / Retrofit_existing_min for the new fossil /
# waste generation technologies.
# ---------------------------------------------------------------------------
$baseLoadTechs = @("Coal", "Coal CCS", "Fuel Oil", "Natural Gas SC", "Natural Gas CC", "Natural Gas CC CCS", "Waste")
$r = 27
foreach ($tech in $baseLoadTechs) {
    $ws.Range("A$r").Value = "CIMS.CAN.BC.Electricity.Utility Generation.Base Load"
    $ws.Range("B$r").Value = "Service"
    $ws.Range("C$r").Value = "BC"
    $ws.Range("D$r").Value = "Electricity"
    $ws.Range("E$r").Value = "Base Load"
    $ws.Range("F$r").Value = $tech
    $ws.Range("G$r").Value = "Retrofit_existing_min"
    $ws.Range("L$r").Value = "%"
    $ws.Range("N$r").Value = 0.5
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Formula = "=O$r"
    $ws.Range("Q$r").Formula = "=P$r"
    $ws.Range("R$r").Formula = "=Q$r"
    $ws.Range("S$r").Formula = "=R$r"
    $ws.Range("T$r").Formula = "=S$r"
    $ws.Range("U$r").Formula = "=T$r"
    $ws.Range("V$r").Formula = "=U$r"
    $ws.Range("W$r").Formula = "=V$r"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# New rows 34-38: Shoulder Load / Retrofit_existing_min.
# ---------------------------------------------------------------------------
$shoulderLoadTechs = @("Coal", "Fuel Oil", "Natural Gas SC", "Natural Gas CC", "Natural Gas CC CCS")
$r = 34
foreach ($tech in $shoulderLoadTechs) {
    $ws.Range("A$r").Value = "CIMS.CAN.BC.Electricity.Utility Generation.Shoulder Load"
    $ws.Range("B$r").Value = "Service"
    $ws.Range("C$r").Value = "BC"
    $ws.Range("D$r").Value = "Electricity"
    $ws.Range("E$r").Value = "Shoulder Load"
    $ws.Range("F$r").Value = $tech
    $ws.Range("G$r").Value = "Retrofit_existing_min"
    $ws.Range("L$r").Value = "%"
    $ws.Range("R$r").Value = 0.5
    $ws.Range("S$r").Value = 1
    $ws.Range("T$r").Formula = "=S$r"
    $ws.Range("U$r").Formula = "=T$r"
    $ws.Range("V$r").Formula = "=U$r"
    $ws.Range("W$r").Formula = "=V$r"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# New rows 39-40: Peak Load / Retrofit_existing_min.
# ---------------------------------------------------------------------------
$peakLoadTechs = @("Natural Gas SC", "Natural Gas CC")
$r = 39
foreach ($tech in $peakLoadTechs) {
    $ws.Range("A$r").Value = "CIMS.CAN.BC.Electricity.Utility Generation.Peak Load"
    $ws.Range("B$r").Value = "Service"
    $ws.Range("C$r").Value = "BC"
    $ws.Range("D$r").Value = "Electricity"
    $ws.Range("E$r").Value = "Peak Load"
    $ws.Range("F$r").Value = $tech
    $ws.Range("G$r").Value = "Retrofit_existing_min"
    $ws.Range("L$r").Value = "%"
    $ws.Range("R$r").Value = 0.5
    $ws.Range("S$r").Value = 1
    $ws.Range("T$r").Formula = "=S$r"
    $ws.Range("U$r").Formula = "=T$r"
    $ws.Range("V$r").Formula = "=U$r"
    $ws.Range("W$r").Formula = "=V$r"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# View tidy-up: select the full used range from A1, no frozen/scrolled cell.
# ---------------------------------------------------------------------------
$ws.Range("A1:X40").Select()
